# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change B11 ("R40") to the text "1" while keeping its original cell
# style/format (General number format, same borders/fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")
$scratch = $ws.Range("ZZ1")

# Stash the cell's current formatting in an out-of-the-way scratch cell so
# we can re-apply it after the value write below (writing a numeric-looking
# string like "1" would otherwise make Excel re-classify the cell as a
# number and drop its existing style).
$target.Copy($scratch)

# Force the cell to Text so assigning "1" is stored as a string (matching
# the original cell, which already held a string) instead of being
# auto-converted to the number 1.
$target.NumberFormat = "@"
$target.Value = "1"

# Restore the original formatting (General number format etc.) on top of
# the new value.
$scratch.Copy()
$target.PasteSpecial(-4122) # xlPasteFormats

# Clean up the scratch cell.
$scratch.Clear()
